# Auto-generated edit script: updates Leve profit-calc sheets with refreshed
# market-price data (per scheduled runner), matching the target commit diff.
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")

# Row 53
$ws.Range("H53").Value = 8063.625
$ws.Range("I53").Value = 1321.2222
$ws.Range("J53").Value = 16732.428
$ws.Range("K53").Value = 1321.2222
$ws.Range("L53").Value = 16732.428
$ws.Range("M53").Value = -684.2221999999999
$ws.Range("N53").Value = -18006.428

# Row 55
$ws.Range("H55").Value = 536.46155
$ws.Range("I55").Value = 548.375
$ws.Range("J55").Value = 517.4
$ws.Range("K55").Value = 548.375
$ws.Range("L55").Value = 517.4
$ws.Range("M55").Value = -334.375
$ws.Range("N55").Value = -945.4

# Row 70
$ws.Range("H70").Value = 2023.5294
$ws.Range("J70").Value = 2500
$ws.Range("L70").Value = 7500
$ws.Range("N70").Value = -8040

# Row 73
$ws.Range("H73").Value = 2023.5294
$ws.Range("J73").Value = 2500
$ws.Range("L73").Value = 7500
$ws.Range("N73").Value = -9372

# Row 74
$ws.Range("H74").Value = 2404.7058
$ws.Range("I74").Value = 1725.3334
$ws.Range("K74").Value = 1725.3334
$ws.Range("M74").Value = -789.3334

# Row 77
$ws.Range("H77").Value = 2404.7058
$ws.Range("I77").Value = 1725.3334
$ws.Range("K77").Value = 8626.666999999999
$ws.Range("M77").Value = -3946.666999999999

# Row 131
$ws.Range("H131").Value = 500569.8
$ws.Range("I131").Value = 500569.8
$ws.Range("K131").Value = 1501709.4
$ws.Range("M131").Value = -1496669.4

# Row 137
$ws.Range("H137").Value = 11180.956
$ws.Range("I137").Value = 1363.4286
$ws.Range("K137").Value = 4090.2858
$ws.Range("M137").Value = -1540.2858

# Row 138
$ws.Range("H138").Value = 2274.3164
$ws.Range("I138").Value = 1940.5555
$ws.Range("J138").Value = 2468.1128
$ws.Range("K138").Value = 5821.666499999999
$ws.Range("L138").Value = 7404.3384
$ws.Range("M138").Value = -681.6664999999994
$ws.Range("N138").Value = -17684.3384

# Row 141
$ws.Range("H141").Value = 3316.1333
$ws.Range("I141").Value = 3410.1428
$ws.Range("K141").Value = 10230.4284
$ws.Range("M141").Value = -5050.428400000001

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")

# Row 132
$ws.Range("H132").Value = 6692979
$ws.Range("I132").Value = 8949.571
$ws.Range("J132").Value = 12541505
$ws.Range("K132").Value = 26848.713
$ws.Range("L132").Value = 37624515
$ws.Range("M132").Value = -24318.713
$ws.Range("N132").Value = -37629575

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")

# Row 23
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()

# Row 105
$ws.Range("H105").Value = 58825016
$ws.Range("I105").Value = 71430104
$ws.Range("J105").Value = 1299.6666
$ws.Range("K105").Value = 71430104
$ws.Range("L105").Value = 1299.6666
$ws.Range("M105").Value = -71428357
$ws.Range("N105").Value = -4793.6666

# Row 134
$ws.Range("H134").Value = 40101
$ws.Range("I134").Value = 47362.293
$ws.Range("K134").Value = 142086.879
$ws.Range("M134").Value = -139551.879

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")

# Row 31
$ws.Range("H31").Value = 9216.784
$ws.Range("I31").Value = 881.8095
$ws.Range("J31").Value = 20156.438
$ws.Range("K31").Value = 881.8095
$ws.Range("L31").Value = 20156.438
$ws.Range("M31").Value = -586.8095
$ws.Range("N31").Value = -20746.438

# Row 34
$ws.Range("H34").Value = 9216.784
$ws.Range("I34").Value = 881.8095
$ws.Range("J34").Value = 20156.438
$ws.Range("K34").Value = 881.8095
$ws.Range("L34").Value = 20156.438
$ws.Range("M34").Value = -679.8095
$ws.Range("N34").Value = -20560.438

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")

# Row 11
$ws.Range("H11").Value = 473.9565
$ws.Range("I11").Value = 123.22222
$ws.Range("J11").Value = 1736.6
$ws.Range("K11").Value = 369.66666
$ws.Range("L11").Value = 5209.799999999999
$ws.Range("M11").Value = -229.66666
$ws.Range("N11").Value = -5489.799999999999

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")

# Row 13
$ws.Range("H13").Value = 481.63635
$ws.Range("I13").Value = 150
$ws.Range("K13").Value = 150
$ws.Range("M13").Value = -11

# Row 80
$ws.Range("H80").Value = 15995.25
$ws.Range("I80").Value = 13289.728
$ws.Range("K80").Value = 13289.728
$ws.Range("M80").Value = -12291.728

# Row 83
$ws.Range("H83").Value = 15995.25
$ws.Range("I83").Value = 13289.728
$ws.Range("K83").Value = 66448.64
$ws.Range("M83").Value = -61456.64

# Row 113
$ws.Range("H113").Value = 3027.2727
$ws.Range("I113").Value = 2538.4
$ws.Range("J113").Value = 3434.6667
$ws.Range("K113").Value = 2538.4
$ws.Range("L113").Value = 3434.6667
$ws.Range("M113").Value = -368.4000000000001
$ws.Range("N113").Value = -7774.6667

# Row 122
$ws.Range("H122").Value = 2616988.2
$ws.Range("I122").Value = 2834321
$ws.Range("J122").Value = 8995
$ws.Range("K122").Value = 8502963
$ws.Range("L122").Value = 26985
$ws.Range("M122").Value = -8500513
$ws.Range("N122").Value = -31885

# Row 123
$ws.Range("H123").Value = 49999
$ws.Range("J123").Value = 49999
$ws.Range("L123").Value = 49999
$ws.Range("N123").Value = -54899

# Row 132
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")

# Row 7
$ws.Range("H7").Value = 2076596.5
$ws.Range("I7").Value = 2386786
$ws.Range("K7").Value = 2386786
$ws.Range("M7").Value = -2386674

# Row 40
$ws.Range("H40").Value = 1734960.9
$ws.Range("I40").Value = 3607.5789
$ws.Range("K40").Value = 3607.5789
$ws.Range("M40").Value = -3471.5789

# Row 68
$ws.Range("H68").Value = 3751.9
$ws.Range("I68").Value = 3242.25
$ws.Range("J68").Value = 4091.6667
$ws.Range("K68").Value = 3242.25
$ws.Range("L68").Value = 4091.6667
$ws.Range("M68").Value = -2493.25
$ws.Range("N68").Value = -5589.6667

# Row 71
$ws.Range("H71").Value = 3751.9
$ws.Range("I71").Value = 3242.25
$ws.Range("J71").Value = 4091.6667
$ws.Range("K71").Value = 16211.25
$ws.Range("L71").Value = 20458.3335
$ws.Range("M71").Value = -12467.25
$ws.Range("N71").Value = -27946.3335

# Row 82
$ws.Range("H82").Value = 2468.5908
$ws.Range("I82").Value = 3088.0908
$ws.Range("J82").Value = 1849.091
$ws.Range("K82").Value = 3088.0908
$ws.Range("L82").Value = 1849.091
$ws.Range("M82").Value = -2727.0908
$ws.Range("N82").Value = -2571.091

# Row 85
$ws.Range("H85").Value = 2468.5908
$ws.Range("I85").Value = 3088.0908
$ws.Range("J85").Value = 1849.091
$ws.Range("K85").Value = 3088.0908
$ws.Range("L85").Value = 1849.091
$ws.Range("M85").Value = -1840.0908
$ws.Range("N85").Value = -4345.091

# Row 122
$ws.Range("H122").Value = 32107160
$ws.Range("I122").Value = 54818404
$ws.Range("K122").Value = 164455212
$ws.Range("M122").Value = -164452762

# Row 126
$ws.Range("H126").Value = 2076596.5
$ws.Range("I126").Value = 2386786
$ws.Range("K126").Value = 7160358
$ws.Range("M126").Value = -7157888

# Row 132
$ws.Range("H132").Value = 1035712.3
$ws.Range("I132").Value = 5173.636
$ws.Range("K132").Value = 15520.908
$ws.Range("M132").Value = -12990.908

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")

# Row 11
$ws.Range("H11").Value = 4677701.5
$ws.Range("I11").Value = 9333668
$ws.Range("J11").Value = 21735
$ws.Range("K11").Value = 9333668
$ws.Range("L11").Value = 21735
$ws.Range("M11").Value = -9333526
$ws.Range("N11").Value = -22019

# Row 13
$ws.Range("H13").Value = 7988.375
$ws.Range("J13").Value = 21000
$ws.Range("L13").Value = 21000
$ws.Range("N13").Value = -21280

# Row 126
$ws.Range("H126").Value = 7148028.5
$ws.Range("I126").Value = 5950.1665
$ws.Range("K126").Value = 17850.4995
$ws.Range("M126").Value = -15380.4995
